$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 5236
$ws.Range("F5").Value = 29
$ws.Range("F8").Value = 586
$ws.Range("F9").Value = 543
$ws.Range("F12").Value = 1445
$ws.Range("F13").Value = 4178
$ws.Range("F14").Value = 431
$ws.Range("F16").Value = 153
$ws.Range("F18").Value = 3239
$ws.Range("F19").Value = 160
$ws.Range("F20").Value = 1069
$ws.Range("F21").Value = 98
$ws.Range("F24").Value = 105
$ws.Range("F25").Value = 32
$ws.Range("F27").Value = 68
$ws.Range("F28").Value = 293
$ws.Range("F31").Value = 14
$ws.Range("F32").Value = 17
$ws.Range("F33").Value = 17
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 5236
$ws.Range("F6").Value = 29
$ws.Range("F9").Value = 586
$ws.Range("F10").Value = 543
$ws.Range("F13").Value = 1445
$ws.Range("F14").Value = 4178
$ws.Range("F15").Value = 431
$ws.Range("F17").Value = 153
$ws.Range("F19").Value = 3239
$ws.Range("F20").Value = 160
$ws.Range("F21").Value = 1069
$ws.Range("F22").Value = 98
$ws.Range("F25").Value = 105
$ws.Range("F26").Value = 32
$ws.Range("F28").Value = 68
$ws.Range("F29").Value = 293
$ws.Range("F32").Value = 14
$ws.Range("F33").Value = 17
$ws.Range("F34").Value = 17
